# Applies the resume.docx edits described by the commit diff:
#   1. Append " Band-8" to the IBM CDL line.
#   2. Note that Java is the strongest skill ("比较精通").
#   3. Trim the Spark/HDFS bullet (drop "但是没有做深入的研究").
#   4. Tighten the kubernetes bullet wording.
#   5. Tighten the machine-learning bullet wording.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

Replace-Text `
    "IBM 中国软件开发中心（CDL）" `
    "IBM 中国软件开发中心（CDL） Band-8"

Replace-Text `
    "java用的时间最长，其他的像go，nodejs, python, shell, scala, javascript, c++在不同的阶段使用过。" `
    "java用的时间最长，比较精通，其他的像go，nodejs, python, shell, scala, javascript, c++在不同的阶段使用过。"

Replace-Text `
    "大数据主要使用过Spark，HDFS等，但是没有做深入的研究。" `
    "大数据主要使用过Spark，HDFS等。"

Replace-Text `
    "容器化，主要是对kubernetes有一些研究，读了些相关代码，根据我们自己的项目需要写过一些定制化的存储插件。" `
    "容器化，主要研究过kubernetes，读相关代码，根据自己的项目需要写过定制化的存储插件。"

Replace-Text `
    "机器学习，目前这个组主要是为数据科学家提供开发平台，不涉及到具体的使用场景，根据自己的兴趣学习中。" `
    "机器学习，目前所在工作组主要是为数据科学家提供开发平台，不涉及到具体的使用场景。"
